$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cxcl12"
$ws.Cells.Item(2,3).Value = "Ackr3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3.0
$ws.Cells.Item(2,6).Value = 1.0
$ws.Cells.Item(2,7).Value = 198.977211
$ws.Cells.Item(2,8).Value = 596.9316329999999
$ws.Cells.Item(2,9).Value = 0.440791350614085
$ws.Cells.Item(2,10).Value = 0.4407913506140851
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 21.90542466666666
$ws.Cells.Item(2,14).Value = 65.716274
$ws.Cells.Item(2,15).Value = 0.4103613011498748
$ws.Cells.Item(2,16).Value = 0.4103613011498748
$ws.Cells.Item(2,17).Value = 4358.680305943937
$ws.Cells.Item(2,18).Value = 39228.12275349544
$ws.Cells.Item(2,19).Value = 0.1808837121736066
$ws.Cells.Item(2,20).Value = 0.1808837121736066

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cxcl12"
$ws.Cells.Item(3,3).Value = "Ackr3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3.0
$ws.Cells.Item(3,6).Value = 1.0
$ws.Cells.Item(3,7).Value = 198.977211
$ws.Cells.Item(3,8).Value = 596.9316329999999
$ws.Cells.Item(3,9).Value = 0.440791350614085
$ws.Cells.Item(3,10).Value = 0.4407913506140851
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 23.19964933333333
$ws.Cells.Item(3,14).Value = 69.598948
$ws.Cells.Item(3,15).Value = 0.4346064242769223
$ws.Cells.Item(3,16).Value = 0.4346064242769223
$ws.Cells.Item(3,17).Value = 4616.201520524675
$ws.Cells.Item(3,18).Value = 41545.81368472207
$ws.Cells.Item(3,19).Value = 0.1915707527425826
$ws.Cells.Item(3,20).Value = 0.1915707527425827

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cxcl12"
$ws.Cells.Item(4,3).Value = "Ackr3"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3.0
$ws.Cells.Item(4,6).Value = 1.0
$ws.Cells.Item(4,7).Value = 198.977211
$ws.Cells.Item(4,8).Value = 596.9316329999999
$ws.Cells.Item(4,9).Value = 0.440791350614085
$ws.Cells.Item(4,10).Value = 0.4407913506140851
$ws.Cells.Item(4,11).Value = 3.0
$ws.Cells.Item(4,12).Value = 1.0
$ws.Cells.Item(4,13).Value = 8.275750666666667
$ws.Cells.Item(4,14).Value = 24.827252
$ws.Cells.Item(4,15).Value = 0.1550322745732029
$ws.Cells.Item(4,16).Value = 0.1550322745732029
$ws.Cells.Item(4,17).Value = 1646.685786584724
$ws.Cells.Item(4,18).Value = 14820.17207926251
$ws.Cells.Item(4,19).Value = 0.06833688569789577
$ws.Cells.Item(4,20).Value = 0.06833688569789577

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Cxcl12"
$ws.Cells.Item(5,3).Value = "Ackr3"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3.0
$ws.Cells.Item(5,6).Value = 1.0
$ws.Cells.Item(5,7).Value = 196.7746736666667
$ws.Cells.Item(5,8).Value = 590.324021
$ws.Cells.Item(5,9).Value = 0.4359121013721307
$ws.Cells.Item(5,10).Value = 0.4359121013721308
$ws.Cells.Item(5,11).Value = 3.0
$ws.Cells.Item(5,12).Value = 1.0
$ws.Cells.Item(5,13).Value = 21.90542466666666
$ws.Cells.Item(5,14).Value = 65.716274
$ws.Cells.Item(5,15).Value = 0.4103613011498748
$ws.Cells.Item(5,16).Value = 0.4103613011498748
$ws.Cells.Item(5,17).Value = 4310.432790313083
$ws.Cells.Item(5,18).Value = 38793.89511281776
$ws.Cells.Item(5,19).Value = 0.1788814571060437
$ws.Cells.Item(5,20).Value = 0.1788814571060437

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cxcl12"
$ws.Cells.Item(6,3).Value = "Ackr3"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3.0
$ws.Cells.Item(6,6).Value = 1.0
$ws.Cells.Item(6,7).Value = 196.7746736666667
$ws.Cells.Item(6,8).Value = 590.324021
$ws.Cells.Item(6,9).Value = 0.4359121013721307
$ws.Cells.Item(6,10).Value = 0.4359121013721308
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 23.19964933333333
$ws.Cells.Item(6,14).Value = 69.598948
$ws.Cells.Item(6,15).Value = 0.4346064242769223
$ws.Cells.Item(6,16).Value = 0.4346064242769223
$ws.Cells.Item(6,17).Value = 4565.103426747767
$ws.Cells.Item(6,18).Value = 41085.9308407299
$ws.Cells.Item(6,19).Value = 0.189450199676381
$ws.Cells.Item(6,20).Value = 0.189450199676381

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cxcl12"
$ws.Cells.Item(7,3).Value = "Ackr3"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 196.7746736666667
$ws.Cells.Item(7,8).Value = 590.324021
$ws.Cells.Item(7,9).Value = 0.4359121013721307
$ws.Cells.Item(7,10).Value = 0.4359121013721308
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 8.275750666666667
$ws.Cells.Item(7,14).Value = 24.827252
$ws.Cells.Item(7,15).Value = 0.1550322745732029
$ws.Cells.Item(7,16).Value = 0.1550322745732029
$ws.Cells.Item(7,17).Value = 1628.458136780033
$ws.Cells.Item(7,18).Value = 14656.12323102029
$ws.Cells.Item(7,19).Value = 0.06758044458970601
$ws.Cells.Item(7,20).Value = 0.06758044458970602

$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Cxcl12"
$ws.Cells.Item(8,3).Value = "Ackr3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 2.0
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.158424
$ws.Cells.Item(8,8).Value = 0.475272
$ws.Cells.Item(8,9).Value = 0.0003509544061791369
$ws.Cells.Item(8,10).Value = 0.0003509544061791369
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 21.90542466666666
$ws.Cells.Item(8,14).Value = 65.716274
$ws.Cells.Item(8,15).Value = 0.4103613011498748
$ws.Cells.Item(8,16).Value = 0.4103613011498748
$ws.Cells.Item(8,17).Value = 3.470344997392
$ws.Cells.Item(8,18).Value = 31.233104976528
$ws.Cells.Item(8,19).Value = 0.0001440181067639523
$ws.Cells.Item(8,20).Value = 0.0001440181067639523

$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Cxcl12"
$ws.Cells.Item(9,3).Value = "Ackr3"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 2.0
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.158424
$ws.Cells.Item(9,8).Value = 0.475272
$ws.Cells.Item(9,9).Value = 0.0003509544061791369
$ws.Cells.Item(9,10).Value = 0.0003509544061791369
$ws.Cells.Item(9,11).Value = 3.0
$ws.Cells.Item(9,12).Value = 1.0
$ws.Cells.Item(9,13).Value = 23.19964933333333
$ws.Cells.Item(9,14).Value = 69.598948
$ws.Cells.Item(9,15).Value = 0.4346064242769223
$ws.Cells.Item(9,16).Value = 0.4346064242769223
$ws.Cells.Item(9,17).Value = 3.675381245984
$ws.Cells.Item(9,18).Value = 33.078431213856
$ws.Cells.Item(9,19).Value = 0.0001525270395537453
$ws.Cells.Item(9,20).Value = 0.0001525270395537453

$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Cxcl12"
$ws.Cells.Item(10,3).Value = "Ackr3"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 2.0
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.158424
$ws.Cells.Item(10,8).Value = 0.475272
$ws.Cells.Item(10,9).Value = 0.0003509544061791369
$ws.Cells.Item(10,10).Value = 0.0003509544061791369
$ws.Cells.Item(10,11).Value = 3.0
$ws.Cells.Item(10,12).Value = 1.0
$ws.Cells.Item(10,13).Value = 8.275750666666667
$ws.Cells.Item(10,14).Value = 24.827252
$ws.Cells.Item(10,15).Value = 0.1550322745732029
$ws.Cells.Item(10,16).Value = 0.1550322745732029
$ws.Cells.Item(10,17).Value = 1.311077523616
$ws.Cells.Item(10,18).Value = 11.799697712544
$ws.Cells.Item(10,19).Value = 0.00005440925986143931
$ws.Cells.Item(10,20).Value = 0.00005440925986143931

$ws.Cells.Item(11,1).Value = "sCs"
$ws.Cells.Item(11,2).Value = "Cxcl12"
$ws.Cells.Item(11,3).Value = "Ackr3"
$ws.Cells.Item(11,4).Value = "ECs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 55.49875533333334
$ws.Cells.Item(11,8).Value = 166.496266
$ws.Cells.Item(11,9).Value = 0.122945593607605
$ws.Cells.Item(11,10).Value = 0.122945593607605
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 21.90542466666666
$ws.Cells.Item(11,14).Value = 65.716274
$ws.Cells.Item(11,15).Value = 0.4103613011498748
$ws.Cells.Item(11,16).Value = 0.4103613011498748
$ws.Cells.Item(11,17).Value = 1215.723804048098
$ws.Cells.Item(11,18).Value = 10941.51423643289
$ws.Cells.Item(11,19).Value = 0.0504521137634605
$ws.Cells.Item(11,20).Value = 0.0504521137634605

$ws.Cells.Item(12,1).Value = "sCs"
$ws.Cells.Item(12,2).Value = "Cxcl12"
$ws.Cells.Item(12,3).Value = "Ackr3"
$ws.Cells.Item(12,4).Value = "FAPs"
$ws.Cells.Item(12,5).Value = 3.0
$ws.Cells.Item(12,6).Value = 1.0
$ws.Cells.Item(12,7).Value = 55.49875533333334
$ws.Cells.Item(12,8).Value = 166.496266
$ws.Cells.Item(12,9).Value = 0.122945593607605
$ws.Cells.Item(12,10).Value = 0.122945593607605
$ws.Cells.Item(12,11).Value = 3.0
$ws.Cells.Item(12,12).Value = 1.0
$ws.Cells.Item(12,13).Value = 23.19964933333333
$ws.Cells.Item(12,14).Value = 69.598948
$ws.Cells.Item(12,15).Value = 0.4346064242769223
$ws.Cells.Item(12,16).Value = 0.4346064242769223
$ws.Cells.Item(12,17).Value = 1287.551662169796
$ws.Cells.Item(12,18).Value = 11587.96495952817
$ws.Cells.Item(12,19).Value = 0.05343294481840483
$ws.Cells.Item(12,20).Value = 0.05343294481840483

$ws.Cells.Item(13,1).Value = "sCs"
$ws.Cells.Item(13,2).Value = "Cxcl12"
$ws.Cells.Item(13,3).Value = "Ackr3"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3.0
$ws.Cells.Item(13,6).Value = 1.0
$ws.Cells.Item(13,7).Value = 55.49875533333334
$ws.Cells.Item(13,8).Value = 166.496266
$ws.Cells.Item(13,9).Value = 0.122945593607605
$ws.Cells.Item(13,10).Value = 0.122945593607605
$ws.Cells.Item(13,11).Value = 3.0
$ws.Cells.Item(13,12).Value = 1.0
$ws.Cells.Item(13,13).Value = 8.275750666666667
$ws.Cells.Item(13,14).Value = 24.827252
$ws.Cells.Item(13,15).Value = 0.1550322745732029
$ws.Cells.Item(13,16).Value = 0.1550322745732029
$ws.Cells.Item(13,17).Value = 459.2938614490037
$ws.Cells.Item(13,18).Value = 4133.644753041032
$ws.Cells.Item(13,19).Value = 0.01906053502573963
$ws.Cells.Item(13,20).Value = 0.01906053502573963

